# Weekly refresh of the Esparragos dataset: rotate the data rows (2-15)
# so that each row takes on the values that used to belong to the row
# 7 positions below it (wrapping around within rows 2-15), for the
# columns D (Fecha) and I..Q (Calidad..Kg o Unidades).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 15
$rowCount = $lastRow - $firstRow + 1
$shift = 7

# Columns whose values move together with each record.
$cols = @("D", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the current ("before") values for every relevant cell so that
# writes to one row don't affect the source data read for another row.
$original = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $original[$r] = $rowVals
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $r + $shift
    if ($srcRow -gt $lastRow) {
        $srcRow = $srcRow - $rowCount
    }
    $srcVals = $original[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
